$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-03-13", "admin", "ali haroon", "01:06:03", "01:06:20", "cleaner"),
    @("2025-03-13", "Ahmed", "Ahmed nawaz", "01:19:04", "01:19:20", "manager"),
    @("2025-03-13", "b3tablocker", "bukhari", "01:19:10", "01:19:16", "trainer")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        if ($c -eq 1) {
            # Column A holds date-like text (e.g. "2025-03-13"); prefix with an
            # apostrophe so Excel keeps it as literal text instead of coercing
            # it into a date serial number, then reset the style so no extra
            # "text" number-format styling gets attached to the cell.
            $cell.Value = "'" + $values[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$c - 1]
        }
    }
}
